$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "297.51"
Set-TextValue $ws.Range("E2") "1.28%"
Set-TextValue $ws.Range("D3") "41.78"
Set-TextValue $ws.Range("E3") "2.96%"
Set-TextValue $ws.Range("D4") "5.019"
Set-TextValue $ws.Range("E4") "-0.12%"
Set-TextValue $ws.Range("D5") "0.07525"
Set-TextValue $ws.Range("E5") "2.65%"
Set-TextValue $ws.Range("D6") "1.602"
Set-TextValue $ws.Range("E6") "4.30%"
Set-TextValue $ws.Range("E7") "-1.12%"
Set-TextValue $ws.Range("E8") "1.74%"
Set-TextValue $ws.Range("E9") "1.42%"
Set-TextValue $ws.Range("D10") "0.1825"
Set-TextValue $ws.Range("E10") "4.50%"
Set-TextValue $ws.Range("D11") "0.08980"
Set-TextValue $ws.Range("E11") "2.93%"
Set-TextValue $ws.Range("D12") "0.04088"
Set-TextValue $ws.Range("E12") "-5.98%"
Set-TextValue $ws.Range("E13") "-0.51%"
Set-TextValue $ws.Range("B14") "BitForexToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001284"
Set-TextValue $ws.Range("E14") "0.97%"
Set-TextValue $ws.Range("B15") "TigerCash"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.005925"
Set-TextValue $ws.Range("E15") "-0.76%"
Set-TextValue $ws.Range("B16") "LEO"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.340"
Set-TextValue $ws.Range("E16") "0.02%"
Set-TextValue $ws.Range("B17") "GateToken"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "4.376"
Set-TextValue $ws.Range("E17") "2.09%"
Set-TextValue $ws.Range("B18") "BitpandaEcosystemToken"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D18") "0.3328"
Set-TextValue $ws.Range("E18") "1.18%"
Set-TextValue $ws.Range("B19") "MCDex"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D19") "8.301"
Set-TextValue $ws.Range("E19") "4.14%"
Set-TextValue $ws.Range("B20") "ProBitToken"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D20") "0.1372"
Set-TextValue $ws.Range("E20") "-1.36%"
Set-TextValue $ws.Range("B21") "ZBToken"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D21") "0.3222"
Set-TextValue $ws.Range("E21") "17.48%"
Set-TextValue $ws.Range("D22") "0.04093"
Set-TextValue $ws.Range("E22") "3.95%"
Set-TextValue $ws.Range("E23") "0.43%"
Set-TextValue $ws.Range("E24") "7.39%"
Set-TextValue $ws.Range("D25") "0.0001302"
Set-TextValue $ws.Range("E25") "1.64%"
Set-TextValue $ws.Range("D38") "0.02405"
Set-TextValue $ws.Range("E38") "4.07%"
Set-TextValue $ws.Range("D39") "0.05205"
Set-TextValue $ws.Range("D40") "0.006307"
Set-TextValue $ws.Range("E40") "1.98%"
Set-TextValue $ws.Range("D41") "0.007814"
Set-TextValue $ws.Range("E41") "-0.54%"
Set-TextValue $ws.Range("D42") "0.1326"
Set-TextValue $ws.Range("E42") "2.96%"
Set-TextValue $ws.Range("D43") "0.007403"
Set-TextValue $ws.Range("E43") "0.74%"
Set-TextValue $ws.Range("D44") "0.007126"
Set-TextValue $ws.Range("E44") "-1.67%"
Set-TextValue $ws.Range("D45") "0.3253"
Set-TextValue $ws.Range("E45") "1.72%"
Set-TextValue $ws.Range("D46") "0.00006589"
Set-TextValue $ws.Range("E46") "4.68%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.04%"
Set-TextValue $ws.Range("D48") "0.04543"
Set-TextValue $ws.Range("E48") "27.37%"
Set-TextValue $ws.Range("D49") "0.004204"
Set-TextValue $ws.Range("E49") "0.06%"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.04%"
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.04%"
